$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 84) with the next stream's stats,
# carrying forward the same formatting as the previous data row.
$ws.Range("A83:D83").Copy()
$ws.Range("A84:D84").PasteSpecial(-4122)

$ws.Range("A84").Value = 82.0
$ws.Range("B84").Value = 247.0
$ws.Range("C84").Value = 333.0
$ws.Range("D84").Value = 105.0
